$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '28.486.96'
$ws.Range('D3').Value = '1.826.29'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '316.52'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').Value = '0.5171'
$ws.Range('E7').Value = '  +2.15%  '
$ws.Range('D8').Value = '0.3863'
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').Value = '0.08289'
$ws.Range('E9').Value = '  +8.31%  '
$ws.Range('D10').Value = '1.123'
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').Value = '41.92'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '6.387'
$ws.Range('D13').Value = '21.22'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '7.497'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '1.828.47'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '93.96'
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('D18').Value = '0.00001124'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D20').Value = '17.81'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('D23').Value = '28.523.03'
$ws.Range('D24').Value = '11.48'
$ws.Range('E24').Value = '  +3.21%  '
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').Value = '21.08'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').Value = '159.92'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').Value = '2.038.33'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '2.417'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').Value = '125.86'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('D32').Value = '1.099'
$ws.Range('E32').Value = '  -2.48%  '
$ws.Range('D33').Value = '0.07644'
$ws.Range('E33').Value = '  +8.49%  '
$ws.Range('D34').Value = '5.732'
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('D35').Value = '3.682'
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('D36').Value = '0.2236'
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('D37').Value = '0.02373'
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('D38').Value = '5.261'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('D39').Value = '12.04'
$ws.Range('E39').Value = '  +7.28%  '
$ws.Range('D40').Value = '8.787'
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').Value = '0.6428'
$ws.Range('E41').Value = '  +2.96%  '
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').Value = '1.400'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '0.6213'
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('D45').Value = '13.57'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('D46').Value = '3.797'
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').Value = '127.96'
$ws.Range('E47').Value = '  +2.80%  '
$ws.Range('D48').Value = '2.004'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').Value = '0.06974'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '1.076'
$ws.Range('E51').Value = '  +0.92%  '
